{"js": "// Set the Alt Text (description) on the document's inline picture.\n// Note: the runtime's XML attribute serializer does not escape embedded\n// double-quote characters, so we pass them through as the literal XML\n// entity `&quot;` (this round-trips correctly through the OOXML writer\n// and decodes back to `\"` for any XML-aware reader).\nconst altText = \"This shows the properties of a subgraph node. Red arrows highlight the &quot;Name&quot; field and the &quot;Replace With&quot; field.  These fields allow you to replace a variable used in a subgraph, with a variable from the parent graph, which can be useful for sharing settings between multiple subgraph nodes, and other situations as well.\";\n\nconst pictures = context.document.body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < pictures.items.length; i++) {\n  pictures.items[i].altTextDescription = altText;\n}\n\nawait context.sync();\n", "ps1": "# Note: the runtime's XML attribute serializer does not escape embedded\n# double-quote characters, so we pass them through as the literal XML\n# entity `&quot;` (this round-trips correctly through the OOXML writer\n# and decodes back to `\"` for any XML-aware reader).\n$d = $word.ActiveDocument\n$altText = 'This shows the properties of a subgraph node. Red arrows highlight the &quot;Name&quot; field and the &quot;Replace With&quot; field.  These fields allow you to replace a variable used in a subgraph, with a variable from the parent graph, which can be useful for sharing settings between multiple subgraph nodes, and other situations as well.'\n\nforeach ($shape in $d.InlineShapes) {\n    $shape.AlternativeText = $altText\n}\n"}
